# Update the "Scoring Card" estimates table for the revised paper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "Sensitive to moments used for estimation ?"
$ws.Range("C2").Value = "Yes"
$ws.Range("E2").Value = "No"

# Row 3 - "Sensitive to assumed inflation process?"
$ws.Range("E3").Value = "No"

# Row 4 - "Sensitive to two-step or joint estimate?"
$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = "No"

# Row 5 - "Sensitive to the type of agents?"
$ws.Range("D5").Value = "Yes"

# Move the active selection to A8, matching the saved cursor position.
$ws.Range("A8").Select()
